$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 40, shifting existing rows 40-43 down to 41-44
$ws.Rows("40").Insert()

# Populate the newly inserted row 40 with the new weekly price entry
$ws.Range("A40").Value2 = 7
$ws.Range("B40").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C40").Value2 = "Ñuble"
$ws.Range("D40").Value2 = 45008
$ws.Range("E40").Value2 = 16
$ws.Range("F40").Value2 = "Fruta"
$ws.Range("G40").Value2 = 100101
$ws.Range("H40").Value2 = "Berries"
$ws.Range("I40").Value2 = 100101001
$ws.Range("J40").Value2 = "Arándano (blue)"
$ws.Range("K40").Value2 = "Sin especificar"
$ws.Range("L40").Value2 = "Primera"
$ws.Range("M40").Value2 = 30
$ws.Range("N40").Value2 = 4000
$ws.Range("O40").Value2 = 4000
$ws.Range("P40").Value2 = 4000
$ws.Range("Q40").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R40").Value2 = "Provincia de Diguillín"
$ws.Range("S40").Value2 = 2000
$ws.Range("T40").Value2 = 2

# Ensure the date cell keeps the same date style/number format as the other date cells (column D)
$ws.Range("D40").NumberFormat = $ws.Range("D41").NumberFormat
